$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation (leading apostrophe) so numeric-looking
# strings like "42.512.72" or "0.0910" are not coerced to numbers,
# matching the original inlineStr text cells.
$ws.Range('D2').Value = "'42.512.72"
$ws.Range('E2').Value = "'  -2.54%  "
$ws.Range('D3').Value = "'2.223.76"
$ws.Range('E3').Value = "'  -2.04%  "
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = "'  +0.18%  "
$ws.Range('D5').Value = "'110.32"
$ws.Range('E5').Value = "'  -6.98%  "
$ws.Range('D6').Value = "'288.84"
$ws.Range('E6').Value = "'  +8.17%  "
$ws.Range('D7').Value = "'0.623"
$ws.Range('E7').Value = "'  -2.95%  "
$ws.Range('E8').Value = "'  -0.36%  "
$ws.Range('D9').Value = "'0.599"
$ws.Range('E9').Value = "'  -3.41%  "
$ws.Range('D10').Value = "'43.44"
$ws.Range('E10').Value = "'  -8.23%  "
$ws.Range('D11').Value = "'0.0910"
$ws.Range('E11').Value = "'  -3.45%  "
$ws.Range('D12').Value = "'54.25"
$ws.Range('E12').Value = "'  +0.85%  "
$ws.Range('D13').Value = "'8.64"
$ws.Range('E13').Value = "'  -8.76%  "
$ws.Range('E14').Value = "'  +13.30%  "
$ws.Range('E15').Value = "'  -2.70%  "
$ws.Range('D16').Value = "'14.86"
$ws.Range('E16').Value = "'  -5.32%  "
$ws.Range('D17').Value = "'2.556.69"
$ws.Range('E17').Value = "'  -2.17%  "
$ws.Range('D18').Value = "'2.230.91"
$ws.Range('E18').Value = "'  -1.67%  "
$ws.Range('D19').Value = "'42.351.70"
$ws.Range('E19').Value = "'  -2.95%  "
$ws.Range('D20').Value = "'7.17"
$ws.Range('E20').Value = "'  +3.75%  "
$ws.Range('E21').Value = "'  -4.38%  "
$ws.Range('D22').Value = "'73.06"
$ws.Range('E22').Value = "'  +0.99%  "
$ws.Range('D23').Value = "'3.35"
$ws.Range('E23').Value = "'  +14.77%  "
$ws.Range('E24').Value = "'  -0.24%  "
$ws.Range('D25').Value = "'231.90"
$ws.Range('E25').Value = "'  -0.99%  "
$ws.Range('D26').Value = "'8.98"
$ws.Range('E26').Value = "'  -5.59%  "
$ws.Range('E27').Value = "'  -1.73%  "
$ws.Range('E28').Value = "'  -6.63%  "
$ws.Range('E29').Value = "'  -2.02%  "
$ws.Range('E30').Value = "'  -4.23%  "
$ws.Range('B31').Value = "'Monero"
$ws.Range('C31').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D31').Value = "'173.05"
$ws.Range('E31').Value = "'  -0.83%  "
$ws.Range('B32').Value = "'InjectiveProtocol"
$ws.Range('C32').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('D32').Value = "'37.07"
$ws.Range('E32').Value = "'  -11.78%  "
$ws.Range('D33').Value = "'20.87"
$ws.Range('E33').Value = "'  -2.72%  "
$ws.Range('D34').Value = "'0.0874"
$ws.Range('E34').Value = "'  -4.29%  "
$ws.Range('D35').Value = "'5.60"
$ws.Range('E35').Value = "'  -2.18%  "
$ws.Range('D36').Value = "'5.01"
$ws.Range('E36').Value = "'  +9.21%  "
$ws.Range('E37').Value = "'  -3.08%  "
$ws.Range('D38').Value = "'4.17"
$ws.Range('E38').Value = "'  -2.13%  "
$ws.Range('D39').Value = "'0.0369"
$ws.Range('E39').Value = "'  -3.76%  "
$ws.Range('E40').Value = "'  -3.63%  "
$ws.Range('D41').Value = "'74.51"
$ws.Range('E41').Value = "'  +2.78%  "
$ws.Range('D42').Value = "'2.39"
$ws.Range('E42').Value = "'  -6.33%  "
$ws.Range('E43').Value = "'  -4.25%  "
$ws.Range('E44').Value = "'  -0.19%  "
$ws.Range('D45').Value = "'12.35"
$ws.Range('E45').Value = "'  -11.00%  "
$ws.Range('E46').Value = "'  -6.01%  "
$ws.Range('E47').Value = "'  -6.24%  "
$ws.Range('D48').Value = "'1.75"
$ws.Range('E48').Value = "'  +12.80%  "
$ws.Range('E49').Value = "'  +1.19%  "
$ws.Range('B50').Value = "'FraxShare"
$ws.Range('C50').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('D50').Value = "'8.46"
$ws.Range('E50').Value = "'  -1.29%  "
$ws.Range('B51').Value = "'Aave"
$ws.Range('C51').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('D51').Value = "'101.75"
$ws.Range('E51').Value = "'  -1.22%  "
